$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N ("Late"), pushing the
# existing N..P columns ("Late", heading, "Outstanding") one to the right.
$ws.Columns("N").Insert()

# Match the width Excel assigns the freshly-inserted column (it inherits
# the neighbouring column's raw width of 11, expressed here in the
# "characters" unit COM uses for ColumnWidth).
$ws.Columns("N").ColumnWidth = 10.166666666666666

# Make "Repayment schedule" the active sheet/tab and move its selection.
$ws.Activate()
$ws.Range("R7").Select()
